$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook records weekly Fruta/Hortaliza price observations for
# "Macroferia Regional de Talca" (Granada). The update reshuffles the
# per-row weekly observations (date, quality, volume, prices, origin, etc.)
# across rows 2-15 while keeping the fixed descriptive columns (A-K) intact.

# Row 2 (now holds the data previously on row 14)
$ws.Range("D2").Value = 44326
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1111
$ws.Range("T2").Value = 18

# Row 3 (now holds the data previously on row 8)
$ws.Range("D3").Value = 44354
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

# Row 4 (now holds the data previously on row 15)
$ws.Range("D4").Value = 44355
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Provincia de Limarí'
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18

# Row 5 (now holds the data previously on row 4)
$ws.Range("D5").Value = 44328
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 1111
$ws.Range("T5").Value = 18

# Row 6 (now holds the data previously on row 6)
$ws.Range("D6").Value = 44299
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 15

# Row 7 (now holds the data previously on row 5)
$ws.Range("D7").Value = 44340
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 230
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1111
$ws.Range("T7").Value = 18

# Row 8 (now holds the data previously on row 9)
$ws.Range("D8").Value = 44358
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 18

# Row 9 (now holds the data previously on row 10)
$ws.Range("D9").Value = 44358
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 17000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 944
$ws.Range("T9").Value = 18

# Row 10 (now holds the data previously on row 3)
$ws.Range("D10").Value = 44291
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 800
$ws.Range("T10").Value = 15

# Row 11 (now holds the data previously on row 2)
$ws.Range("D11").Value = 44348
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 1111
$ws.Range("T11").Value = 18

# Row 12 (now holds the data previously on row 13)
$ws.Range("D12").Value = 44316
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 1111
$ws.Range("T12").Value = 18

# Row 13 (now holds the data previously on row 7)
$ws.Range("D13").Value = 44342
$ws.Range("L13").Value = 'Especial'
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 1111
$ws.Range("T13").Value = 18

# Row 14 (now holds the data previously on row 11)
$ws.Range("D14").Value = 44294
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = '$/caja 15 kilos granel'
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 800
$ws.Range("T14").Value = 15

# Row 15 (now holds the data previously on row 12)
$ws.Range("D15").Value = 44319
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 1111
$ws.Range("T15").Value = 18
